$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '58.060.45'
$ws.Range("E2").Value = '  -2.94%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.301.85'
$ws.Range("E3").Value = '  -3.38%  '

# Row 4: TetherUSD
$c = $ws.Range("D4")
$c.Value = "'" + '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

# Row 5: BNB
$c = $ws.Range("D5")
$c.Value = "'" + '535.25'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.77%  '

# Row 6: Solana
$c = $ws.Range("D6")
$c.Value = "'" + '131.02'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.86%  '

# Row 7: USDC
$c = $ws.Range("D7")
$c.Value = "'" + '1.00'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '

# Row 8: XRP
$ws.Range("E8").Value = '  -1.20%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '2.299.48'
$ws.Range("E9").Value = '  -3.42%  '

# Row 10: Dogecoin
$c = $ws.Range("D10")
$c.Value = "'" + '0.0997'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -5.07%  '

# Row 11: Toncoin
$c = $ws.Range("D11")
$c.Value = "'" + '5.43'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -3.20%  '

# Row 12: TRON
$c = $ws.Range("D12")
$c.Value = "'" + '0.149'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.29%  '

# Row 13: Cardano
$ws.Range("E13").Value = '  -3.87%  '

# Row 14: Avalanche
$c = $ws.Range("D14")
$c.Value = "'" + '23.48'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.89%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.707.47'
$ws.Range("E15").Value = '  -3.61%  '

# Row 16: WrappedBTC
$ws.Range("D16").Value = '58.012.02'
$ws.Range("E16").Value = '  -2.93%  '

# Row 17: ShibaInu
$ws.Range("E17").Value = '  -3.88%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.293.46'
$ws.Range("E18").Value = '  -3.66%  '

# Row 19: Chainlink
$c = $ws.Range("D19")
$c.Value = "'" + '10.52'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -5.34%  '

# Row 20: Polkadot
$c = $ws.Range("D20")
$c.Value = "'" + '4.22'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -5.75%  '

# Row 21: BitcoinCash
$c = $ws.Range("D21")
$c.Value = "'" + '313.64'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.11%  '

# Row 22: Uniswap
$c = $ws.Range("D22")
$c.Value = "'" + '6.37'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -4.90%  '

# Row 23: Dai
$ws.Range("E23").Value = '  -0.17%  '

# Row 24: Litecoin
$c = $ws.Range("D24")
$c.Value = "'" + '62.67'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.32%  '

# Row 25: Kaspa
$ws.Range("E25").Value = '  -3.82%  '

# Row 26: Binance-PegBSC-USD
$c = $ws.Range("D26")
$c.Value = "'" + '0.999'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '

# Row 27: InternetComputer(DFINITY)
$c = $ws.Range("D27")
$c.Value = "'" + '8.02'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -4.54%  '

# Row 28: Fetch.AI
$c = $ws.Range("D28")
$c.Value = "'" + '1.30'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -5.07%  '

# Row 29: Monero
$c = $ws.Range("D29")
$c.Value = "'" + '170.57'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.48%  '

# Row 30: PancakeSwap
$ws.Range("E30").Value = '  -5.36%  '

# Row 31: PEPE
$ws.Range("D31").Value = '0.0₃0720'
$ws.Range("E31").Value = '  -5.07%  '

# Row 32: Aptos
$c = $ws.Range("D32")
$c.Value = "'" + '5.78'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -4.58%  '

# Row 33: SuiNetwork
$c = $ws.Range("D33")
$c.Value = "'" + '1.05'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.47%  '

# Row 34: PolygonEcosystemToken
$c = $ws.Range("D34")
$c.Value = "'" + '0.377'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.67%  '

# Row 35: USDe
$c = $ws.Range("D35")
$c.Value = "'" + '0.999'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '

# Row 36: EthereumClassic
$ws.Range("E36").Value = '  -2.31%  '

# Row 37: FirstDigitalUSD
$ws.Range("E37").Value = '  -0.10%  '

# Row 38: ImmutableX
$c = $ws.Range("D38")
$c.Value = "'" + '1.24'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -6.47%  '

# Row 39: NEARProtocol
$c = $ws.Range("D39")
$c.Value = "'" + '3.90'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -5.70%  '

# Row 40: OKB
$c = $ws.Range("D40")
$c.Value = "'" + '38.05'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.47%  '

# Row 41: Stacks
$ws.Range("E41").Value = '  -5.28%  '

# Row 42: Aave
$c = $ws.Range("D42")
$c.Value = "'" + '141.00'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -3.27%  '

# Row 43: Bittensor
$c = $ws.Range("D43")
$c.Value = "'" + '288.67'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -9.45%  '

# Row 44: Filecoin
$ws.Range("E44").Value = '  -3.10%  '

# Row 45: Stellar
$c = $ws.Range("D45")
$c.Value = "'" + '0.0948'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.19%  '

# Row 46: Hedera
$c = $ws.Range("D46")
$c.Value = "'" + '0.0497'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.56%  '

# Row 47: Mantle
$ws.Range("E47").Value = '  -2.94%  '

# Row 48: InjectiveProtocol
$c = $ws.Range("D48")
$c.Value = "'" + '18.14'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -7.89%  '

# Row 49: VeChain
$c = $ws.Range("D49")
$c.Value = "'" + '0.0211'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -3.04%  '

# Row 50: WhiteBITCoin
$ws.Range("E50").Value = '  -1.22%  '

# Row 51: BabyDogeCoin
$ws.Range("D51").Value = '0.0₆0202'
$ws.Range("E51").Value = '  +84.64%  '
